# "fix the merge error" - restore the View (column F) flags that were lost
# in a bad merge on the "Property" sheet, and correct a couple of
# Public/Private flags that slipped along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Rows 68-75 lost their "View" (column F) flag during the merge - restore it.
for ($r = 68; $r -le 75; $r++) {
    $ws.Cells.Item($r, 6).Value = $true
}

# Row 76 (GameID): Public/Private were wrongly kept TRUE and View wrongly FALSE.
$ws.Cells.Item(76, 4).Value = $false   # Private
$ws.Cells.Item(76, 5).Value = $false   # Save
$ws.Cells.Item(76, 6).Value = $true    # View

# Row 77 (GateID): same fix as row 76.
$ws.Cells.Item(77, 4).Value = $false   # Private
$ws.Cells.Item(77, 5).Value = $false   # Save
$ws.Cells.Item(77, 6).Value = $true    # View

# Row 78 (GuildID): Public flag was wrongly left TRUE.
$ws.Cells.Item(78, 3).Value = $false   # Public

# Re-point the selection left over from the merge to where the fix was made.
$ws.Activate()
$ws.Range("C78").Select()
